$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column L values, one per existing row, keyed off the matching
# column K cell so the new cells inherit the exact same formatting/style
# that the rest of the row already uses.
$pairs = @(
    @{src = "K4";  dst = "L4";  val = 2023},
    @{src = "K5";  dst = "L5";  val = 1.6430457248453274},
    @{src = "K6";  dst = "L6";  val = 0.41181606829870221},
    @{src = "K7";  dst = "L7";  val = 0.94796963217320562},
    @{src = "K8";  dst = "L8";  val = 0.72306112208737106},
    @{src = "K9";  dst = "L9";  val = 2.1802539701246277},
    @{src = "K10"; dst = "L10"; val = 0.63651150401750112},
    @{src = "K11"; dst = "L11"; val = 0.97994201681774651},
    @{src = "K12"; dst = "L12"; val = 2.2469385026996971},
    @{src = "K13"; dst = "L13"; val = 4.1686356866605365},
    @{src = "K14"; dst = "L14"; val = 0.3304193846038968}
)

foreach ($p in $pairs) {
    # Copy formatting (style) from the column K cell onto the new
    # column L cell, then set the new cell's own value.
    $ws.Range($p.src).Copy() | Out-Null
    $ws.Range($p.dst).PasteSpecial(-4122)
    $ws.Range($p.dst).Value = $p.val
}

# Rows 2 and 3 grew slightly taller in the edited workbook.
$ws.Rows(2).RowHeight = 13.5
$ws.Rows(3).RowHeight = 13.5

Write-Output "done"
